$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '27.292.27'
$ws.Cells.Item(2, 5).Value = '  -2.37%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.708.64'
$ws.Cells.Item(3, 5).Value = '  -1.72%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.004'
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''223.58'
$ws.Cells.Item(5, 5).Value = '  -2.59%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''0.5309'
$ws.Cells.Item(6, 5).Value = '  -2.59%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''1.003'
$ws.Cells.Item(7, 5).Value = '  -0.02%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.2656'
$ws.Cells.Item(8, 5).Value = '  -4.49%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.06594'
$ws.Cells.Item(9, 5).Value = '  -1.82%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''20.89'
$ws.Cells.Item(10, 5).Value = '  -4.09%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.07647'
$ws.Cells.Item(11, 5).Value = '  -1.78%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''4.582'
$ws.Cells.Item(12, 5).Value = '  -2.83%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.721.98'
$ws.Cells.Item(13, 5).Value = '  -0.90%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.944.82'
$ws.Cells.Item(14, 5).Value = '  -1.73%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''0.5739'
$ws.Cells.Item(15, 5).Value = '  -4.70%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '0.0₅8195'
$ws.Cells.Item(16, 5).Value = '  -2.64%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''67.58'
$ws.Cells.Item(17, 5).Value = '  -2.95%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '27.297.19'
$ws.Cells.Item(18, 5).Value = '  -2.37%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''216.42'
$ws.Cells.Item(19, 5).Value = '  -3.37%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''4.677'
$ws.Cells.Item(21, 5).Value = '  -3.09%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''10.41'
$ws.Cells.Item(22, 5).Value = '  -4.98%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''5.975'
$ws.Cells.Item(23, 5).Value = '  -4.40%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.08%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''142.10'
$ws.Cells.Item(25, 5).Value = '  -3.17%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +5.71%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''0.1218'
$ws.Cells.Item(27, 5).Value = '  -2.55%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''7.262'
$ws.Cells.Item(28, 5).Value = '  -2.79%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -4.53%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''0.05390'
$ws.Cells.Item(30, 5).Value = '  -4.66%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -2.03%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -5.27%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -3.47%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.98%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''2.881'
$ws.Cells.Item(35, 5).Value = '  +0.71%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -1.09%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''0.9481'
$ws.Cells.Item(37, 5).Value = '  -3.71%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''0.5856'
$ws.Cells.Item(38, 5).Value = '  -1.63%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''0.01635'
$ws.Cells.Item(39, 5).Value = '  -2.56%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''5.862'
$ws.Cells.Item(40, 5).Value = '  -2.67%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Maker'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(41, 4).Value = '1.043.60'
$ws.Cells.Item(41, 5).Value = '  -0.24%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'PaxDollar'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(42, 4).Value = '''1.004'
$ws.Cells.Item(42, 5).Value = '  +0.04%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.8394'
$ws.Cells.Item(43, 5).Value = '  -0.99%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''101.09'
$ws.Cells.Item(44, 5).Value = '  -1.17%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '1.851.06'
$ws.Cells.Item(45, 5).Value = '  -1.74%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '0.0₈115'
$ws.Cells.Item(46, 5).Value = '  +5.30%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -3.44%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''0.4501'
$ws.Cells.Item(48, 5).Value = '  +1.59%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -0.65%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''8.108'
$ws.Cells.Item(50, 5).Value = '  -2.26%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''0.06500'
$ws.Cells.Item(51, 5).Value = '  +9.08%  '
